$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Weekly topics shift up by one row, since week 2 used to span three rows
# (3,4,5) and now only occupies one (row 3): collapse the old three-topic
# write-up into a single shorter title, then shift the rest up. ---
$ws.Range("F3").Value = "The file system; the shell; the terminal"
$ws.Range("F4").Value = "Editing text: Text editors; regular expressions"
$ws.Range("F5").Value = "Your data workbench I: R, RStudio, and Quarto"
$ws.Range("F6").Value = "Your data workbench II: How R thinks; tidy data"
$ws.Range("F7").Value = "Version Control: git and GitHub"
$ws.Range("F8").Value = "No class (Fall break)"
$ws.Range("F9").Value = "Wrangle data: Getting stuff in and out of R"
$ws.Range("F10").Value = "Tabulate data: Grouping, summaries"
$ws.Range("F11").Value = "Look at data: Graphs, ggplot, and the grammar of graphics"
$ws.Range("F12").Value = "Iterate on data: functional programming patterns"
$ws.Range("F13").Value = "Reproducible results: build systems, environments, and packages"

# Row 14 no longer has a topic (the list is one row shorter than before).
$ws.Range("F14").ClearContents()

# Week 2 (row 3) gains content/example/assignment links, matching week 1's
# pattern in row 2.
$ws.Range("H3").Value = "/content/02-content"
$ws.Range("I3").Value = "/example/02-example"
$ws.Range("J3").Value = "/assignment/02-assignment"
# J2 (week 1's assignment cell) uses the bold/black "applyFont" style; match it on J3.
$ws.Range("J3").Font.Color = 0

# Update the saved cursor/selection position to reflect where editing left off.
$ws.Range("J4").Select()
